# Rename the "Order" list values from top/bottom to above/below.
#
# The workbook has a hidden-ish lookup sheet named "." whose column B
# ("order") holds the two allowed values used by the data-validation list
# on the "Layers" sheet's "Order" column (E). We rename "top" -> "above"
# and "bottom" -> "below" there, then update every cell on "Layers" that
# currently holds the old value so nothing is left pointing at stale text.

$wb = $excel.ActiveWorkbook

$wsLookup = $wb.Worksheets.Item(".")
$wsLayers = $wb.Worksheets.Item("Layers")

# 1) Rename the two list values in the lookup sheet.
$wsLookup.Range("B2").Value = "above"
$wsLookup.Range("B3").Value = "below"

# 2) Propagate the rename to every "Order" cell (column E) on the Layers
#    sheet that still references the old wording.
$lastRow = $wsLayers.Cells(1, 1).SpecialCells(11).Row
for ($r = 2; $r -le $lastRow; $r++) {
    $cell = $wsLayers.Cells.Item($r, 5)
    $val = $cell.Value2
    if ($val -eq "top") {
        $cell.Value = "above"
    } elseif ($val -eq "bottom") {
        $cell.Value = "below"
    }
}

# 3) Match the author's final UI state: active cell on the lookup sheet is
#    B3, but the workbook still shows "Layers" as the active tab.
$wsLookup.Range("B3").Select() | Out-Null
$wsLayers.Activate()
